$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Exact "  ARIYAMANGALAM                                      Retailer Name   :  Thulasi-D                          " `
              "  ARIYAMANGALAM                                      Retailer Name   :  ESHWAR MEDICALS                    "

Replace-Exact "  TRICHY-620010                                      Address         :                                     " `
              "  TRICHY-620010                                      Address         :  Thamarai 1st St Ezhil Nagar Tamil  "

Replace-Exact "  PHONE NO         :9944951444                                          Ariyamangalam                      " `
              "  PHONE NO         :9944951444                                          India Near Kumutha Store           "

Replace-Exact "  GSTIN No         :33AAPFD1365C1ZR                                     TRICHY                             " `
              "  GSTIN No         :33AAPFD1365C1ZR                                                                        "

Replace-Exact "  RS PAN No        :AAPFD1365C                       Phone No        :   9943684000                        " `
              "  RS PAN No        :AAPFD1365C                       Phone No        :   8072005857                        "

Replace-Exact "  Salesperson Name :N SANTHOSH                       " `
              "  Salesperson Name :SAKTHIVEL M                      "

Replace-Exact "  Beat Name        :Chemist - Kattur NUTS            GSTIN NO        :   33AABCT6876B1ZF                   " `
              "  Beat Name        :Chemist - Thiruvarambur NUTS     GSTIN NO        :                                     "

Replace-Exact "  HUL STORE ID     :HUL-41A392D-P19230               Time of Billing :   07/09/2023 09:57:35               " `
              "  HUL STORE ID     :HUL-41A392D-P25120               Time of Billing :   14/12/2023 22:26:36               "

Replace-Exact "  Five Thousand Seven Hundred Forty Rupees Only                              " `
              "  Six Hundred Twenty-Nine Rupees Only                                        "

Replace-Exact "  ABC32856     Thulasi-D    Amt : 5740.00" `
              "  ABC54009     ESHWAR MEDICALS    Amt : 629.00"

Write-Host "Done applying replacements"
